$d = $word.ActiveDocument

$replacements = @(
    @("2024-09-21 Saturday", "2024-09-22 Sunday"),
    @("898÷8=", "958÷3="),
    @("121÷5=", "686÷3="),
    @("910÷4=", "266÷9="),
    @("887÷2=", "861÷5="),
    @("407÷2=", "162÷9="),
    @("732÷9=", "828÷2="),
    @("525÷2=", "753÷9="),
    @("336÷2=", "597÷3="),
    @("666÷2=", "421÷4="),
    @("554÷5=", "431÷9="),
    @("878÷9=", "976÷6="),
    @("844÷5=", "431÷2="),
    @("181÷2=", "683÷2="),
    @("632÷6=", "893÷5="),
    @("763÷5=", "108÷5="),
    @("322÷9=", "585÷7="),
    @("202÷4=", "786÷3="),
    @("680÷9=", "661÷2="),
    @("534÷4=", "828÷3="),
    @("668÷8=", "673÷6="),
    @("434÷6=", "626÷8="),
    @("849÷8=", "246÷4="),
    @("972÷7=", "738÷3="),
    @("444÷2=", "964÷4="),
    @("691÷2=", "677÷8=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
